$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2322.3462
$ws.Range("I15").Value = 2322.3462
$ws.Range("K15").Value = 6967.0386
$ws.Range("M15").Value = -6798.0386

$ws.Range("H33").Value = 768.4286
$ws.Range("I33").Value = 371.2857
$ws.Range("J33").Value = 1165.5714
$ws.Range("K33").Value = 371.2857
$ws.Range("L33").Value = 1165.5714
$ws.Range("M33").Value = -142.2857
$ws.Range("N33").Value = -1623.5714

$ws.Range("H53").Value = 465.53845
$ws.Range("I53").Value = 123.2
$ws.Range("K53").Value = 123.2
$ws.Range("M53").Value = 513.8

$ws.Range("H86").Value = 3475.5334
$ws.Range("I86").Value = 2931
$ws.Range("J86").Value = 4292.3335
$ws.Range("K86").Value = 2931
$ws.Range("L86").Value = 4292.3335
$ws.Range("M86").Value = -1808
$ws.Range("N86").Value = -6538.3335

$ws.Range("H89").Value = 3475.5334
$ws.Range("I89").Value = 2931
$ws.Range("J89").Value = 4292.3335
$ws.Range("K89").Value = 14655
$ws.Range("L89").Value = 21461.6675
$ws.Range("M89").Value = -9039
$ws.Range("N89").Value = -32693.6675

$ws.Range("H98").Value = 902.9
$ws.Range("I98").Value = 966.1875
$ws.Range("K98").Value = 966.1875
$ws.Range("M98").Value = 531.8125

$ws.Range("H107").Value = 1685.7
$ws.Range("I107").Value = 2116.7693
$ws.Range("K107").Value = 2116.7693
$ws.Range("M107").Value = -196.7692999999999

$ws.Range("H122").Value = 902.9
$ws.Range("I122").Value = 966.1875
$ws.Range("K122").Value = 2898.5625
$ws.Range("M122").Value = -448.5625

$ws.Range("H132").Value = 12533.019
$ws.Range("I132").Value = 1951.619
$ws.Range("K132").Value = 5854.857
$ws.Range("M132").Value = -3324.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16133249
$ws.Range("I32").Value = 16133249
$ws.Range("K32").Value = 16133249
$ws.Range("M32").Value = -16132962

$ws.Range("H74").Value = 2114.4075
$ws.Range("I74").Value = 2128.9583
$ws.Range("J74").Value = 1998
$ws.Range("K74").Value = 2128.9583
$ws.Range("L74").Value = 1998
$ws.Range("M74").Value = -1254.9583
$ws.Range("N74").Value = -3746

$ws.Range("H77").Value = 2114.4075
$ws.Range("I77").Value = 2128.9583
$ws.Range("J77").Value = 1998
$ws.Range("K77").Value = 10644.7915
$ws.Range("L77").Value = 9990
$ws.Range("M77").Value = -6276.791499999999
$ws.Range("N77").Value = -18726

$ws.Range("H97").Value = 2274
$ws.Range("I97").Value = 1258.6666
$ws.Range("J97").Value = 4304.6665
$ws.Range("K97").Value = 1258.6666
$ws.Range("L97").Value = 4304.6665
$ws.Range("M97").Value = -762.6666
$ws.Range("N97").Value = -5296.6665

$ws.Range("H109").Value = 49999.832
$ws.Range("J109").Value = 49999.832
$ws.Range("L109").Value = 49999.832
$ws.Range("N109").Value = -52773.832

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2632.6538
$ws.Range("I94").Value = 2938.6
$ws.Range("K94").Value = 2938.6
$ws.Range("M94").Value = -2487.6

$ws.Range("H134").Value = 2531.8572
$ws.Range("I134").Value = 2245.7368
$ws.Range("K134").Value = 6737.2104
$ws.Range("M134").Value = -4202.2104

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 140.93333
$ws.Range("I7").Value = 49.166668
$ws.Range("J7").Value = 202.11111
$ws.Range("K7").Value = 49.166668
$ws.Range("L7").Value = 202.11111
$ws.Range("M7").Value = 63.833332
$ws.Range("N7").Value = -428.11111

$ws.Range("H48").Value = 39999.75
$ws.Range("J48").Value = 39999.75
$ws.Range("L48").Value = 39999.75
$ws.Range("N48").Value = -40951.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 2615.6667
$ws.Range("I14").Value = 2615.6667
$ws.Range("K14").Value = 7847.000100000001
$ws.Range("M14").Value = -7674.000100000001

$ws.Range("H136").Value = 2607.6428
$ws.Range("I136").Value = 813.375
$ws.Range("K136").Value = 2440.125
$ws.Range("M136").Value = 2659.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws.Range("H105").Value = 95000
$ws.Range("J105").Value = 95000
$ws.Range("L105").Value = 95000
$ws.Range("N105").Value = -101988

$ws.Range("H107").Value = 681.5
$ws.Range("I107").Value = 681.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 681.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1238.5
$ws.Range("N107").ClearContents()

$ws.Range("H113").Value = 3064.55
$ws.Range("I113").Value = 1830.6875
$ws.Range("K113").Value = 1830.6875
$ws.Range("M113").Value = 339.3125

$ws.Range("H132").Value = 1628.6
$ws.Range("I132").Value = 1373.8636
$ws.Range("K132").Value = 4121.5908
$ws.Range("M132").Value = -1591.5908

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5657.241
$ws.Range("I7").Value = 5420.25
$ws.Range("K7").Value = 5420.25
$ws.Range("M7").Value = -5308.25

$ws.Range("H9").Value = 942
$ws.Range("I9").Value = 403.33334
$ws.Range("J9").Value = 1750
$ws.Range("K9").Value = 403.33334
$ws.Range("L9").Value = 1750
$ws.Range("M9").Value = -179.33334
$ws.Range("N9").Value = -2198

$ws.Range("H58").Value = 4537.2
$ws.Range("I58").Value = 1395.3334
$ws.Range("J58").Value = 9250
$ws.Range("K58").Value = 1395.3334
$ws.Range("L58").Value = 9250
$ws.Range("M58").Value = -1135.3334
$ws.Range("N58").Value = -9770

$ws.Range("H126").Value = 5657.241
$ws.Range("I126").Value = 5420.25
$ws.Range("K126").Value = 16260.75
$ws.Range("M126").Value = -13790.75

$ws.Range("H132").Value = 6739
$ws.Range("I132").Value = 2260.7
$ws.Range("K132").Value = 6782.099999999999
$ws.Range("M132").Value = -4252.099999999999

$ws.Range("H136").Value = 2876.9688
$ws.Range("J136").Value = 2224.5
$ws.Range("L136").Value = 6673.5
$ws.Range("N136").Value = -11773.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 11511.6
$ws.Range("I45").Value = 9784
$ws.Range("K45").Value = 9784
$ws.Range("M45").Value = -9293

$ws.Range("I81").Value = 4279.6665
$ws.Range("K81").Value = 8559.333000000001
$ws.Range("M81").Value = -7498.333000000001

$ws.Range("I84").Value = 4279.6665
$ws.Range("K84").Value = 42796.665
$ws.Range("M84").Value = -37492.665

$ws.Range("H96").Value = 35198.125
$ws.Range("I96").Value = 86357
$ws.Range("J96").Value = 4502.8
$ws.Range("K96").Value = 86357
$ws.Range("L96").Value = 4502.8
$ws.Range("M96").Value = -84984
$ws.Range("N96").Value = -7248.8

$ws.Range("H136").Value = 3178.0715
$ws.Range("I136").Value = 1856.2858
$ws.Range("K136").Value = 5568.857400000001
$ws.Range("M136").Value = -3018.857400000001
